$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Update the "datetimeFigureOut" date placeholders on every slide layout
#    and on the slide master: 7/18/2021 -> 7/21/2021
# ---------------------------------------------------------------------------
$cls = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $cls.Count; $li++) {
    $layout = $cls.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        if ($sh.PlaceholderFormat.Type -eq 16) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "7/18/2021") {
                $tr.Text = "7/21/2021"
            }
        }
    }
}

$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.PlaceholderFormat.Type -eq 16) {
        $tr = $sh.TextFrame.TextRange
        if ($tr.Text -eq "7/18/2021") {
            $tr.Text = "7/21/2021"
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Slide 10 ("מחבר ישר 8" connector): reposition the line
# ---------------------------------------------------------------------------
$s10 = $p.Slides.Item(10)
for ($i = 1; $i -le $s10.Shapes.Count; $i++) {
    $sh = $s10.Shapes.Item($i)
    if ($sh.Name -eq "מחבר ישר 8") {
        $sh.Left = 9.964252068503937
        $sh.Top = 240.16787401574803
    }
}

# ---------------------------------------------------------------------------
# 3) Slide 12 text fix: "הצלחנו" -> "הצלחתי"
# ---------------------------------------------------------------------------
$s12 = $p.Slides.Item(12)
for ($i = 1; $i -le $s12.Shapes.Count; $i++) {
    $sh = $s12.Shapes.Item($i)
    if (-not $sh.HasTextFrame) { continue }
    $tr = $sh.TextFrame.TextRange
    $paras = $tr.Paragraphs()
    for ($pi = 1; $pi -le $paras.Count; $pi++) {
        $para = $tr.Paragraphs($pi, 1)
        if ($para.Text.Length -ge 6) {
            $word = $para.Characters(1, 6)
            if ($word.Text -eq "הצלחנו") {
                $word.Text = "הצלחתי"
            }
        }
    }
}
